$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: Polarverhalten fertiggestellt und getested / 3h / 43400
$ws.Range("B20").Value = "Polarverhalten fertiggestellt und getested"
$ws.Range("D20").Value = "3h"
$ws.Range("F20").Value = 43400
$ws.Range("F20").NumberFormat = $ws.Range("F18").NumberFormat()

# Row 21: Recherche von Collision detection / 2h / 43400
$ws.Range("B21").Value = "Recherche von Collision detection"
$ws.Range("D21").Value = "2h"
$ws.Range("F21").Value = 43400
$ws.Range("F21").NumberFormat = $ws.Range("F18").NumberFormat()

# Update selection to follow the new last row, like Excel would after data entry
$ws.Range("F22").Select()
